# Weekly update: insert the newest day's price record for Berenjena
# (Vega Monumental Concepción) at the top of the data block (row 47),
# pushing all existing data rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 47:112 down to 48:113, inserting a fresh blank row at 47.
$ws.Range("A47:R47").Insert()

# Populate the newly inserted row with the latest reading.
$ws.Range("A47").Value = 11
$ws.Range("B47").Value = "Vega Monumental Concepción"
$ws.Range("C47").Value = "Bíobío"
$ws.Range("D47").Value = 44803
$ws.Range("E47").Value = 8
$ws.Range("F47").Value = 100112001
$ws.Range("G47").Value = "Berenjena"
$ws.Range("H47").Value = "Sin especificar"
$ws.Range("I47").Value = "Primera"
$ws.Range("J47").Value = 110
$ws.Range("K47").Value = 15000
$ws.Range("L47").Value = 16000
$ws.Range("M47").Value = 15545
$ws.Range("N47").Value = "$/caja 60 unidades"
$ws.Range("O47").Value = "Región de Arica y Parinacota"
$ws.Range("P47").Value = 259
$ws.Range("Q47").Value = 60
$ws.Range("R47").Value = "Hortaliza"
